$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (SAIDAS GERAL -> Saidas Geral). Excel keeps the
# _FilterDatabase defined name / autofilter reference in sync automatically.
$ws.Name = "Saidas Geral"

# Move the current selection from X12 to J10 (pane stays frozen at A12).
$ws.Range("J10").Select()

# Apply the "Separador de milhares" (comma / #,##0.00) number format to the
# totals row (9), subtotal row (10) and header row (12) across columns M:W.
$numFmt = "_-* #,##0.00_-;\-* #,##0.00_-;_-* ""-""??_-;_-@_-"
$ws.Range("M9:W9").NumberFormat = $numFmt
$ws.Range("M10:W10").NumberFormat = $numFmt
$ws.Range("M12:W12").NumberFormat = $numFmt
